$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.022.83"
$ws.Range("E2").Value = "  +5.47%  "
$ws.Range("D3").Value = "2.539.63"
$ws.Range("E3").Value = "  +6.52%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'505.60"
$ws.Range("E5").Value = "  +5.65%  "
$ws.Range("D6").Value = "'159.59"
$ws.Range("E6").Value = "  +8.10%  "
$ws.Range("D7").Value = "'0.616"
$ws.Range("E7").Value = "  +23.05%  "
$ws.Range("D8").Value = "'0.993"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "2.580.46"
$ws.Range("E9").Value = "  +8.02%  "
$ws.Range("D10").Value = "'6.25"
$ws.Range("E10").Value = "  +14.47%  "
$ws.Range("E11").Value = "  +7.02%  "
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "2.976.95"
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D15").Value = "58.891.73"
$ws.Range("E15").Value = "  +5.10%  "
$ws.Range("D16").Value = "'22.04"
$ws.Range("E16").Value = "  +8.46%  "
$ws.Range("E17").Value = "  +5.05%  "
$ws.Range("D18").Value = "2.573.32"
$ws.Range("E18").Value = "  +8.14%  "
$ws.Range("E19").Value = "  +4.91%  "
$ws.Range("D20").Value = "'333.97"
$ws.Range("E20").Value = "  +6.12%  "
$ws.Range("D21").Value = "'10.37"
$ws.Range("E21").Value = "  +6.68%  "
$ws.Range("E22").Value = "  +6.74%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "'59.65"
$ws.Range("E24").Value = "  +4.95%  "
$ws.Range("D25").Value = "'0.418"
$ws.Range("E25").Value = "  +5.76%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +6.21%  "
$ws.Range("D27").Value = "'0.990"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "2.637.32"
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("D29").Value = "'7.58"
$ws.Range("E29").Value = "  +4.61%  "
$ws.Range("E30").Value = "  +7.47%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").Value = "'19.48"
$ws.Range("E32").Value = "  +8.14%  "
$ws.Range("D33").Value = "'155.03"
$ws.Range("E33").Value = "  +4.90%  "
$ws.Range("E34").Value = "  +5.81%  "
$ws.Range("E35").Value = "  +9.21%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.20"
$ws.Range("E36").Value = "  +8.68%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'3.94"
$ws.Range("E37").Value = "  +10.06%  "
$ws.Range("D38").Value = "'0.861"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("E39").Value = "  +10.18%  "
$ws.Range("E40").Value = "  +7.13%  "
$ws.Range("D41").Value = "'291.47"
$ws.Range("D42").Value = "'34.71"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("E43").Value = "  +7.28%  "
$ws.Range("E44").Value = "  +7.10%  "
$ws.Range("E45").Value = "  +4.20%  "
$ws.Range("D46").Value = "'0.990"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  +7.48%  "
$ws.Range("D48").Value = "'19.32"
$ws.Range("E48").Value = "  +13.61%  "
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'0.718"
$ws.Range("E50").Value = "  +14.33%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.26"
$ws.Range("E51").Value = "  +0.31%  "
